$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -188.6720833333333
$ws.Range("C2").Value = -999805.1822153517
$ws.Range("D2").Value = 20622.70629166667
$ws.Range("E2").Value = $false
